# Applies two changes to "Portfolio website planning.pptx":
#   1. The "datetimeFigureOut" date field cached on every slide layout
#      (and the slide master) moves from 7/3/20 to 7/5/20.
#   2. On slide 1, the "Tag" schema box ("Rectangle 4") had a row reading
#      "type<TAB>string" — the column named "type" is renamed to
#      "category", splitting that row's single run into two runs:
#      "category" and "<TAB>string".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Update every cached datetimeFigureOut placeholder (7/3/20 -> 7/5/20)
# ---------------------------------------------------------------------
function Update-DateField($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            $full = $tr.Text
            $idx = $full.IndexOf("7/3/20")
            if ($idx -ge 0) {
                $start = $idx + 1
                $sub = $tr.Characters($start, 6)
                $sub.Text = "7/5/20"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateField $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $lyt = $layouts.Item($li)
    Update-DateField $lyt.Shapes
}

# ---------------------------------------------------------------------
# 2) Rename the "type" column to "category" in the Tag schema rectangle
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(1)
$tagShape = $slide.Shapes.Item(2)   # "Rectangle 4"

$tr = $tagShape.TextFrame.TextRange
$full = $tr.Text
$marker = "type" + [char]9 + "string"
$idx = $full.IndexOf($marker)
if ($idx -ge 0) {
    $start = $idx + 1
    $sub = $tr.Characters($start, 4)   # exactly the word "type"
    $sub.Text = "category"
}
